$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_5301_topic_2__ID** " -> "**ID__AFFARS_SUBPART_5301_1__ID**" ---
# The existing text spans two runs: the placeholder text, then a trailing space run.
# Replacing across both in one Find/Execute collapses them into a single run using
# the first run's formatting, which also removes the now-superfluous trailing-space run.
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("**ID__AFFARS_5301_topic_2__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5301_1__ID**", 2)

# --- Update the first paragraph's indentation: w:ind w:left 120 -> 225 (twips => points) ---
# 225 twentieths-of-a-point = 11.25 points
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- Add a paragraph border (pBdr) around the first paragraph, 5pt space on every side ---
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
